# A new price record was added for Jengibre (Vega Modelo de Temuco).
# It slots in right above the current row 233, pushing that row and
# every row after it down by one (233->234, 234->235, ... 322->323).
# We reproduce that by duplicating row 233 (Copy + Insert, which shifts
# everything below it down and leaves the new row 233 as an exact copy
# of the old one) and then updating just the date (column D) of the
# newly inserted row to the new record's date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(233).Copy()
$ws.Rows(233).Insert()

$ws.Range("D233").Value = 45119
